# Update weekly triaged issues
# Rewrites rows 5-54 of the component stats table (Component, Bugs, Features, Total)
# with the refreshed weekly triage numbers / component ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Tree", 10, 0, 10),
    @("DataGrid", 4, 5, 9),
    @("Popover", 6, 2, 8),
    @("Dialog", 5, 2, 7),
    @("Dropdown", 4, 3, 7),
    @("Toolbar", 3, 1, 4),
    @("Virtualizer", 4, 0, 4),
    @("Drawer", 4, 0, 4),
    @("Table", 4, 0, 4),
    @("Tooltip", 3, 1, 4),
    @("TeachingPopover", 3, 0, 3),
    @("Skeleton", 2, 1, 3),
    @("TagPicker", 2, 1, 3),
    @("Portal", 2, 0, 2),
    @("DatePicker", 0, 2, 2),
    @("FluentProvider", 0, 2, 2),
    @("Switch", 2, 0, 2),
    @("Tabs", 2, 0, 2),
    @("Slider", 0, 2, 2),
    @("List", 2, 0, 2),
    @("MessageBar", 2, 0, 2),
    @("Calendar Compat", 2, 0, 2),
    @("Carousel", 1, 0, 1),
    @("Spinner", 0, 1, 1),
    @("Card", 1, 0, 1),
    @("Button", 0, 1, 1),
    @("Input", 0, 1, 1),
    @("InfoLabel", 1, 0, 1),
    @("Accordion", 1, 0, 1),
    @("Popup", 1, 0, 1),
    @("Checkbox", 1, 0, 1),
    @("SplitButton", 1, 0, 1),
    @("AvatarGroup", 1, 0, 1),
    @("Badge", 1, 0, 1),
    @("MenuItem", 1, 0, 1),
    @("Image", 0, 1, 1),
    @("FocusTrapZone", 1, 0, 1),
    @("Label", 1, 0, 1),
    @("Avatar", 1, 0, 1),
    @("DatePickerCompat", 1, 0, 1),
    @("Toast", 1, 0, 1),
    @("Pickers", 0, 0, 0),
    @("Rating", 0, 0, 0),
    @("Tag", 0, 0, 0),
    @("Calendar", 0, 0, 0),
    @("Coachmark", 0, 0, 0),
    @("Keytip", 0, 0, 0),
    @("ColorPicker", 0, 0, 0),
    @("Segment", 0, 0, 0),
    @("SpinButton", 0, 0, 0)
)

$row = 5
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}
